$d = $word.ActiveDocument

# The document ends with an empty paragraph that only contains a manual
# line break (<w:br/>). Add a brand-new paragraph right after it (i.e.
# at the very end of the document body) containing the period total
# line, matching the style used by the other bold headline paragraphs
# in the report ("Dnevni izveštaj za: ...": size 28 half-points / 14pt,
# bold, with spacing before the paragraph).

$text = "UKUPNO ZA PERIOD od 01.06.2025. do 01.07.2025.: 90.000 RSD"

$endRange = $d.Content
$endRange.Collapse(0)               # wdCollapseEnd
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.SpaceBefore = 15           # w:spacing w:before="300" (twentieths of a point)

$newRange = $newPara.Range
$startPos = $newRange.Start
$newRange.InsertAfter($text)

$textRange = $d.Range($startPos, $startPos + $text.Length)
$textRange.Font.Size = 14           # w:sz w:val="28" (half-points)
$textRange.Font.Bold = $true        # w:b w:val="on"
